# Apply the case-locations-and-outbreaks update:
#  - Insert a new "Hallam" row (alphabetically before the "Lakes Entrance" rows)
#  - Insert new "Moorabbin" and "Mordialloc" rows (alphabetically before "Narre Warren"/after "Melbourne")
#  - Append a new "Wonthaggi" row at the end of the table
#
# We work from the bottom of the sheet upward so that earlier inserts do not
# shift the row numbers that later (higher-up) inserts rely on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append "Wonthaggi" row at the end (after the last existing row, 28) ---
$ws.Rows("29:29").Insert()
$ws.Range("A29").Value = "Wonthaggi"
$ws.Range("B29").Value = "Wonthaggi Plaza Shopping Centre, 2 Biggs Drive, Wonthaggi, VIC 3995"
$ws.Range("C29").Value = "28/12/20 1:30pm - 2:30pm"
$ws.Range("D29").Value = "Kmart - shopped for 15 mins"

# --- 2. Insert "Mordialloc" row before old row 22 (Narre Warren) ---
$ws.Rows("22:22").Insert()
$ws.Range("A22").Value = "Mordialloc"
$ws.Range("B22").Value = "Woodlands Golf Club, 109 White Street, Mordialloc, VIC 3195"
$ws.Range("C22").Value = "28/12/20 12:00pm - 6:00pm"
$ws.Range("D22").Value = "Case attended course"

# --- 3. Insert "Moorabbin" row before old row 22 (Narre Warren), i.e. before the Mordialloc row just added ---
$ws.Rows("22:22").Insert()
$ws.Range("A22").Value = "Moorabbin"
$ws.Range("B22").Value = "Costco Moorabbin, 8 Chifley Drive, Moorabbin Airport, VIC 3194"
$ws.Range("C22").Value = "30/12/20 10:45am - 12:15pm"
$ws.Range("D22").Value = "Case shopped in store"

# --- 4. Insert "Hallam" row before old row 9 (Lakes Entrance) ---
$ws.Rows("9:9").Insert()
$ws.Range("A9").Value = "Hallam"
$ws.Range("B9").Value = "Coles Hallam, 2 Princes Domain Drive, Hallam, VIC 3803"
$ws.Range("C9").Value = "30/12/20 6:15am - 6:30am"
$ws.Range("D9").Value = "Case shopped in store"
